$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z3").Value = 42
$x = $ws.Range("Z3").Value
Write-Host ("Z3 is " + $x)
Write-Host $ws.Range("Z3").Text
